# fix(docgen): add alias keys for associe fields and contract date to improve
# template matching.
#
# For each of the three "database" sheets (Societes, Associes, Contrats) the
# last data row (row 9) gets a handful of its numeric / date cells rewritten
# as literal text (so template placeholders resolve against either the
# numeric or the textual alias), and a brand-new row 10 is appended with a
# fresh record.
#
# NOTE: a handful of the literal strings we need to write are "numeric
# looking" (e.g. "9", "100000", "0661545632") or "date looking"
# (e.g. "10/01/2029"). Excel's COM Value setter auto-converts those into
# numbers/dates, same as typing them into the grid would. We force them to
# stay literal text the same way a human would in the UI: a leading
# apostrophe (quote-prefix).

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force literal text storage, even for numeric-/date-looking strings,
    # by using Excel's quote-prefix convention.
    $cell.Value = "'" + $text
}

# ---------------------------------------------------------------------
# Sheet "Societes"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Societes")

# -- row 9: alias a couple of fields as text --
Set-TextValue $ws.Cells.Item(9, 1) "9"
Set-TextValue $ws.Cells.Item(9, 5) "2025-10-28 00:00:00"

# -- row 10: new record (FIRST BUILD) --
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "FIRST BUILD"
$ws.Cells.Item(10, 3).Value = "SARL AU"
Set-TextValue $ws.Cells.Item(10, 4) "00000125465610225"
$ws.Cells.Item(10, 5).NumberFormat = $ws.Cells.Item(9, 5).NumberFormat
$ws.Cells.Item(10, 5).Value = 45958
$ws.Cells.Item(10, 6).Value = "10 000"
Set-TextValue $ws.Cells.Item(10, 7) "100"
$ws.Cells.Item(10, 8).Value = "86 Ha CASABLANCA"
$ws.Cells.Item(10, 9).Value = "Casablanca"

# ---------------------------------------------------------------------
# Sheet "Associes"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Associes")

# -- row 9: alias a couple of fields as text --
Set-TextValue $ws.Cells.Item(9, 1) "9"
Set-TextValue $ws.Cells.Item(9, 2) "9"
Set-TextValue $ws.Cells.Item(9, 14) "1000"
Set-TextValue $ws.Cells.Item(9, 15) "100000"
Set-TextValue $ws.Cells.Item(9, 16) "1"

# -- row 10: new record (SAWAB Alim) --
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = 10
$ws.Cells.Item(10, 3).Value = "M."
$ws.Cells.Item(10, 4).Value = "SAWAB"
$ws.Cells.Item(10, 5).Value = "Alim"
$ws.Cells.Item(10, 6).Value = "Marocaine"
$ws.Cells.Item(10, 7).Value = "G54887"
Set-TextValue $ws.Cells.Item(10, 8) "10/01/2029"
$ws.Cells.Item(10, 9).NumberFormat = $ws.Cells.Item(5, 9).NumberFormat
$ws.Cells.Item(10, 9).Value = 44177
$ws.Cells.Item(10, 10).Value = "DOUR SALIM"
$ws.Cells.Item(10, 11).Value = "Soualem HAi Ajax Willam"
Set-TextValue $ws.Cells.Item(10, 12) "0661545632"
$ws.Cells.Item(10, 13).Value = "Zaimm@gmail.com"
$ws.Cells.Item(10, 14).Value = 1000
$ws.Cells.Item(10, 15).Value = 100000
$ws.Cells.Item(10, 16).Value = 1
$ws.Cells.Item(10, 17).Value = "Associé Gérant"

# ---------------------------------------------------------------------
# Sheet "Contrats"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contrats")

# -- row 9: alias a couple of fields as text --
Set-TextValue $ws.Cells.Item(9, 1) "9"
Set-TextValue $ws.Cells.Item(9, 2) "9"
Set-TextValue $ws.Cells.Item(9, 3) "2025-10-30 00:00:00"
Set-TextValue $ws.Cells.Item(9, 7) "2025-10-30 00:00:00"
Set-TextValue $ws.Cells.Item(9, 8) "2026-10-30 00:00:00"

# -- row 10: new record --
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = 10
$ws.Cells.Item(10, 3).NumberFormat = $ws.Cells.Item(9, 3).NumberFormat
$ws.Cells.Item(10, 3).Value = 45958
Set-TextValue $ws.Cells.Item(10, 4) "15"
$ws.Cells.Item(10, 5).Value = 800
$ws.Cells.Item(10, 6).Value = 1500
$ws.Cells.Item(10, 7).NumberFormat = $ws.Cells.Item(9, 7).NumberFormat
$ws.Cells.Item(10, 7).Value = 45958
$ws.Cells.Item(10, 8).NumberFormat = $ws.Cells.Item(9, 8).NumberFormat
$ws.Cells.Item(10, 8).Value = 46415
